# Updated question 1 and convex function
#
# Replace the placeholder paragraph "<Add rest of Question 1 Instructions
# here>" under "Question 1" with the real instructions describing the
# objective function, split across five paragraphs (two of which are
# blank spacer paragraphs), matching formatting copied from a web source
# (Arial, #222222, 9pt).

$d = $word.ActiveDocument

# Locate the placeholder paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*<Add*rest of*Question 1 Instructions here>*") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Add rest of Question 1 Instructions here' placeholder paragraph."
}

$r = $target.Range

# Build the replacement body: 5 paragraphs (the original paragraph's mark
# is consumed/replaced along with its content since $target.Range already
# spans through the paragraph end).
$body = '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr><w:r><w:rPr><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>Your objective function is:</w:t></w:r></w:p>' +
        '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>' +
        '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr>' +
          '<w:r><w:rPr><w:rFonts w:ascii="arial;sans-serif" w:hAnsi="arial;sans-serif"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:tab/></w:r>' +
          '<w:r><w:rPr><w:rFonts w:ascii="arial;sans-serif" w:hAnsi="arial;sans-serif"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:tab/></w:r>' +
          '<w:r><w:rPr><w:rFonts w:ascii="arial;sans-serif" w:hAnsi="arial;sans-serif"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:tab/></w:r>' +
          '<w:r><w:rPr><w:rFonts w:ascii="arial;sans-serif" w:hAnsi="arial;sans-serif"/><w:b/><w:bCs/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>y = x^5 + 27x^4 + 37x^3 - 339x^2 - 326x + 600</w:t></w:r>' +
        '</w:p>' +
        '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr></w:p>' +
        '<w:p><w:pPr><w:spacing w:after="0"/></w:pPr>' +
          '<w:proofErr w:type="gramStart"/>' +
          '<w:r><w:rPr><w:rFonts w:ascii="arial;sans-serif" w:hAnsi="arial;sans-serif"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t>which</w:t></w:r>' +
          '<w:proofErr w:type="gramEnd"/>' +
          '<w:r><w:rPr><w:rFonts w:ascii="arial;sans-serif" w:hAnsi="arial;sans-serif"/><w:color w:val="222222"/><w:sz w:val="18"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> contains more than one min points. Your task is to try out different step-sizes and starting points and eventually converge to the global minimal point. </w:t></w:r>' +
        '</w:p>'

$frag = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
          '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
              '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                '<w:body>' + $body + '</w:body>' +
              '</w:document>' +
            '</pkg:xmlData>' +
          '</pkg:part>' +
        '</pkg:package>'

$r.InsertXML($frag)
